# Updates the cryptos price-list sheet (Coin/Link/Price/Volume(1h))
# with the latest scrape values from the GitHub Actions job.
# D-column "Price" cells hold dot-grouped numeric-looking text (e.g. "69.355.22")
# that must stay plain text, so each one is forced to Text format before the
# write and reset to the Normal style right after so no stray formatting lingers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.355.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.482.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.89%  "
$ws.Range("E12").Value = "  -3.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.051.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "604.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.473.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.31%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.523.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.87%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.36%  "
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.979"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "105.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.45%  "
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.92%  "
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.06%  "
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.40%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.613.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "509.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.34%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.11%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.391"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.71%  "
$ws.Range("E42").Value = "  -7.27%  "
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("E44").Value = "  -3.65%  "
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("E47").Value = "  -4.76%  "
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.86%  "
